# Apply the Dec 16 2023 "cryptos" data refresh (GitHub Actions job).
# Only cell VALUES change (Coin / Link / Price / Volume(1h) columns);
# no rows/columns are inserted or removed and no formatting changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "42.282.61" },
    @{ Cell = "E2"; Value = "  -0.87%  " },
    @{ Cell = "D3"; Value = "2.242.40" },
    @{ Cell = "E3"; Value = "  -1.12%  " },
    @{ Cell = "E4"; Value = "  -0.06%  " },
    @{ Cell = "D5"; Value = "246.44" },
    @{ Cell = "E5"; Value = "  -1.21%  " },
    @{ Cell = "E6"; Value = "  -1.95%  " },
    @{ Cell = "D7"; Value = "74.67" },
    @{ Cell = "E7"; Value = "  -2.60%  " },
    @{ Cell = "E8"; Value = "  +0.13%  " },
    @{ Cell = "D9"; Value = "0.617" },
    @{ Cell = "E9"; Value = "  -3.27%  " },
    @{ Cell = "D10"; Value = "42.41" },
    @{ Cell = "E10"; Value = "  +6.68%  " },
    @{ Cell = "D11"; Value = "0.0943" },
    @{ Cell = "E11"; Value = "  -2.48%  " },
    @{ Cell = "D12"; Value = "7.15" },
    @{ Cell = "E12"; Value = "  -1.70%  " },
    @{ Cell = "E13"; Value = "  -2.86%  " },
    @{ Cell = "D14"; Value = "14.45" },
    @{ Cell = "E14"; Value = "  -3.28%  " },
    @{ Cell = "D15"; Value = "0.851" },
    @{ Cell = "E15"; Value = "  -1.06%  " },
    @{ Cell = "D16"; Value = "2.239.72" },
    @{ Cell = "E16"; Value = "  -1.05%  " },
    @{ Cell = "D17"; Value = "42.070.56" },
    @{ Cell = "E17"; Value = "  -1.16%  " },
    @{ Cell = "D18"; Value = "0.0₃0985" },
    @{ Cell = "E18"; Value = "  -0.22%  " },
    @{ Cell = "D19"; Value = "6.15" },
    @{ Cell = "E19"; Value = "  -0.18%  " },
    @{ Cell = "D20"; Value = "72.09" },
    @{ Cell = "E20"; Value = "  +0.09%  " },
    @{ Cell = "D21"; Value = "231.43" },
    @{ Cell = "E21"; Value = "  -1.57%  " },
    @{ Cell = "E22"; Value = "  +3.78%  " },
    @{ Cell = "D23"; Value = "8.87" },
    @{ Cell = "E23"; Value = "  +38.88%  " },
    @{ Cell = "E24"; Value = "  +0.07%  " },
    @{ Cell = "D25"; Value = "11.29" },
    @{ Cell = "E25"; Value = "  +0.44%  " },
    @{ Cell = "E26"; Value = "  -4.29%  " },
    @{ Cell = "D27"; Value = "2.30" },
    @{ Cell = "E27"; Value = "  -2.58%  " },
    @{ Cell = "D28"; Value = "169.04" },
    @{ Cell = "E28"; Value = "  +1.04%  " },
    @{ Cell = "E29"; Value = "  -3.92%  " },
    @{ Cell = "D30"; Value = "20.67" },
    @{ Cell = "E30"; Value = "  -0.78%  " },
    @{ Cell = "D31"; Value = "0.0818" },
    @{ Cell = "E31"; Value = "  -3.97%  " },
    @{ Cell = "E32"; Value = "  -2.34%  " },
    @{ Cell = "D33"; Value = "30.44" },
    @{ Cell = "E33"; Value = "  -0.90%  " },
    @{ Cell = "E34"; Value = "  -1.60%  " },
    @{ Cell = "D35"; Value = "5.22" },
    @{ Cell = "E35"; Value = "  +11.12%  " },
    @{ Cell = "D36"; Value = "4.44" },
    @{ Cell = "E36"; Value = "  -2.45%  " },
    @{ Cell = "E37"; Value = "  +3.10%  " },
    @{ Cell = "D38"; Value = "13.55" },
    @{ Cell = "E38"; Value = "  -1.14%  " },
    @{ Cell = "E39"; Value = "  -3.17%  " },
    @{ Cell = "D40"; Value = "5.77" },
    @{ Cell = "E40"; Value = "  -1.12%  " },
    @{ Cell = "D41"; Value = "62.00" },
    @{ Cell = "E41"; Value = "  +1.78%  " },
    @{ Cell = "E42"; Value = "  -1.58%  " },
    @{ Cell = "D43"; Value = "106.56" },
    @{ Cell = "E43"; Value = "  -2.30%  " },
    @{ Cell = "E44"; Value = "  +2.46%  " },
    @{ Cell = "D45"; Value = "8.66" },
    @{ Cell = "E45"; Value = "  -1.91%  " },
    @{ Cell = "D46"; Value = "0.997" },
    @{ Cell = "E46"; Value = "  -0.25%  " },
    @{ Cell = "D47"; Value = "1.12" },
    @{ Cell = "E47"; Value = "  -2.70%  " },
    @{ Cell = "D48"; Value = "4.27" },
    @{ Cell = "E48"; Value = "  -7.50%  " },
    @{ Cell = "E49"; Value = "  -0.17%  " },
    @{ Cell = "D50"; Value = "2.27" },
    @{ Cell = "E50"; Value = "  +1.62%  " },
    @{ Cell = "B51"; Value = "BitTorrent-New" },
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/w4MqH_Xe8+bittorrent-new-btt" },
    @{ Cell = "D51"; Value = "0.0₃0142" },
    @{ Cell = "E51"; Value = "  +11.75%  " }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $isNumericLooking = $u.Value -match '^-?\d+(\.\d+)?$'
    if ($isNumericLooking) {
        # These "Price" figures are stored as literal text in the
        # workbook (e.g. "7.15", "62.00", "106.56") so trailing
        # zeros / exact formatting survive. Assigning the bare
        # numeric-looking string to .Value would make Excel parse
        # it into a real number and drop that formatting, so we
        # use the classic leading-apostrophe "force text" input,
        # then reset the style so no quote-prefix / number-format
        # is stamped onto the cell (it stays the default style,
        # matching the original unstyled cell).
        $range.Value = "'" + $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
